$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.472.57'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.164.91'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.99%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("E6").Value = '  +1.09%  '

$ws.Range("E7").Value = '  +1.86%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  +1.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0853'
$ws.Range("D10").ClearFormats()

$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.13'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.485.15'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.07'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.815'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("E16").Value = '  +0.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.165.75'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.461.42'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.68%  '

$ws.Range("E19").Value = '  +1.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.90'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.70'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.67%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.32'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.11%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.70'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("E28").Value = '  +2.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.87'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("E30").Value = '  +0.40%  '

$ws.Range("E31").Value = '  +6.48%  '

$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.64'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.91%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.73'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("B35").Value = 'THORChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.07'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0622'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.73%  '

$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.60'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '103.21'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.24%  '

$ws.Range("E41").Value = '  +0.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.84'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.98%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.523.55'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("E44").Value = '  +4.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.12'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.93%  '

$ws.Range("E46").Value = '  +1.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0927'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.68%  '

$ws.Range("E48").Value = '  +4.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.368.50'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.98%  '

$ws.Range("E51").Value = '  -0.65%  '
